# New crime data collected - update CompStat_1 sheet figures for the
# 47th Precinct weekly report (week ending 9/8/2024).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Header text: volume/number and reporting week dates (shared strings)
# ---------------------------------------------------------------------
$ws.Cells.Item(8,1).Value = "Volume 31   Number  36"
$ws.Cells.Item(9,3).Value = "Report Covering the Week  9/2/2024  Through  9/8/2024"

# ---------------------------------------------------------------------
# 2) Column E got narrower (bestFit recalculated against the new,
#    shorter figures) - match the other 6.168446-char columns.
# ---------------------------------------------------------------------
$ws.Columns.Item(5).ColumnWidth = 5.43

# ---------------------------------------------------------------------
# 3) Row 22 (Transit): C/D/E flip from numbers to the "no data" text
#    markers already used elsewhere in the sheet ("0" / "***.*"), while
#    F/G pick up new counts. Copy a same-styled text cell first so the
#    style index (14) and shared-string type survive, then set text.
# ---------------------------------------------------------------------
$ws.Cells.Item(14,3).Copy($ws.Cells.Item(22,3))
$ws.Cells.Item(22,3).Value = "0"
$ws.Cells.Item(14,4).Copy($ws.Cells.Item(22,4))
$ws.Cells.Item(22,4).Value = "0"
$ws.Cells.Item(22,5).Copy($ws.Cells.Item(22,5))
$ws.Cells.Item(14,5).Copy($ws.Cells.Item(22,5))
$ws.Cells.Item(22,5).Value = "***.*"

# ---------------------------------------------------------------------
# 4) Row 28 (Shooting Vic.): C28 flips from the "0" text marker to a
#    real numeric count. Copy a same-styled numeric cell first so the
#    style index (15) and numeric type survive, then set the value.
# ---------------------------------------------------------------------
$ws.Cells.Item(27,3).Copy($ws.Cells.Item(28,3))

# ---------------------------------------------------------------------
# 5) Bulk numeric updates, rows 15-30.
# ---------------------------------------------------------------------
$ws.Cells.Item(15,3).Value = 2
$ws.Cells.Item(15,6).Value = 6
$ws.Cells.Item(15,7).Value = 2
$ws.Cells.Item(15,8).Value = 200
$ws.Cells.Item(15,9).Value = 30
$ws.Cells.Item(15,11).Value = 0
$ws.Cells.Item(15,12).Value = -18.918918918918
$ws.Cells.Item(15,13).Value = 3.448275862068
$ws.Cells.Item(15,14).Value = -42.307692307692

$ws.Cells.Item(16,3).Value = 9
$ws.Cells.Item(16,4).Value = 9
$ws.Cells.Item(16,5).Value = 0
$ws.Cells.Item(16,6).Value = 41
$ws.Cells.Item(16,7).Value = 37
$ws.Cells.Item(16,8).Value = 10.810810810810
$ws.Cells.Item(16,9).Value = 360
$ws.Cells.Item(16,10).Value = 323
$ws.Cells.Item(16,11).Value = 11.455108359133
$ws.Cells.Item(16,12).Value = 17.263843648208
$ws.Cells.Item(16,13).Value = 27.208480565371
$ws.Cells.Item(16,14).Value = -62.224554039874

$ws.Cells.Item(17,3).Value = 11
$ws.Cells.Item(17,4).Value = 22
$ws.Cells.Item(17,5).Value = -50
$ws.Cells.Item(17,6).Value = 60
$ws.Cells.Item(17,7).Value = 78
$ws.Cells.Item(17,8).Value = -23.076923076923
$ws.Cells.Item(17,9).Value = 541
$ws.Cells.Item(17,10).Value = 577
$ws.Cells.Item(17,11).Value = -6.239168110918
$ws.Cells.Item(17,12).Value = 4.038461538461
$ws.Cells.Item(17,13).Value = 85.910652920962
$ws.Cells.Item(17,14).Value = -13.162118780096

$ws.Cells.Item(18,3).Value = 7
$ws.Cells.Item(18,4).Value = 2
$ws.Cells.Item(18,5).Value = 250
$ws.Cells.Item(18,6).Value = 17
$ws.Cells.Item(18,7).Value = 19
$ws.Cells.Item(18,8).Value = -10.526315789473
$ws.Cells.Item(18,9).Value = 193
$ws.Cells.Item(18,10).Value = 195
$ws.Cells.Item(18,11).Value = -1.025641025641
$ws.Cells.Item(18,12).Value = -3.5
$ws.Cells.Item(18,13).Value = -15.720524017467
$ws.Cells.Item(18,14).Value = -85.233358837031

$ws.Cells.Item(19,3).Value = 10
$ws.Cells.Item(19,4).Value = 16
$ws.Cells.Item(19,5).Value = -37.5
$ws.Cells.Item(19,6).Value = 66
$ws.Cells.Item(19,7).Value = 55
$ws.Cells.Item(19,8).Value = 20
$ws.Cells.Item(19,9).Value = 621
$ws.Cells.Item(19,10).Value = 532
$ws.Cells.Item(19,11).Value = 16.729323308270
$ws.Cells.Item(19,12).Value = 20.582524271844
$ws.Cells.Item(19,13).Value = 192.924528301887
$ws.Cells.Item(19,14).Value = 64.285714285714

$ws.Cells.Item(20,3).Value = 16
$ws.Cells.Item(20,4).Value = 4
$ws.Cells.Item(20,5).Value = 300
$ws.Cells.Item(20,6).Value = 58
$ws.Cells.Item(20,7).Value = 41
$ws.Cells.Item(20,8).Value = 41.463414634146
$ws.Cells.Item(20,9).Value = 384
$ws.Cells.Item(20,10).Value = 435
$ws.Cells.Item(20,11).Value = -11.724137931034
$ws.Cells.Item(20,12).Value = 23.870967741935
$ws.Cells.Item(20,13).Value = 72.972972972973
$ws.Cells.Item(20,14).Value = -65.560538116591

$ws.Cells.Item(21,3).Value = 55
$ws.Cells.Item(21,4).Value = 53
$ws.Cells.Item(21,5).Value = 3.773584905660
$ws.Cells.Item(21,6).Value = 248
$ws.Cells.Item(21,7).Value = 232
$ws.Cells.Item(21,8).Value = 6.896551724137
$ws.Cells.Item(21,9).Value = 2133
$ws.Cells.Item(21,10).Value = 2102
$ws.Cells.Item(21,11).Value = 1.474785918173
$ws.Cells.Item(21,12).Value = 12.263157894736
$ws.Cells.Item(21,13).Value = 65.992217898832
$ws.Cells.Item(21,14).Value = -52.045863309352

$ws.Cells.Item(22,6).Value = 3
$ws.Cells.Item(22,7).Value = 1
$ws.Cells.Item(22,8).Value = 200

$ws.Cells.Item(23,3).Value = 1
$ws.Cells.Item(23,5).Value = -66.666666666666
$ws.Cells.Item(23,9).Value = 77
$ws.Cells.Item(23,10).Value = 76
$ws.Cells.Item(23,11).Value = 1.315789473684
$ws.Cells.Item(23,12).Value = -4.938271604938
$ws.Cells.Item(23,13).Value = 45.283018867924

$ws.Cells.Item(24,3).Value = 15
$ws.Cells.Item(24,4).Value = 19
$ws.Cells.Item(24,5).Value = -21.052631578947
$ws.Cells.Item(24,6).Value = 82
$ws.Cells.Item(24,7).Value = 104
$ws.Cells.Item(24,8).Value = -21.153846153846
$ws.Cells.Item(24,9).Value = 887
$ws.Cells.Item(24,10).Value = 931
$ws.Cells.Item(24,11).Value = -4.726100966702
$ws.Cells.Item(24,12).Value = -15.200764818355
$ws.Cells.Item(24,13).Value = 73.2421875

$ws.Cells.Item(25,3).Value = 2
$ws.Cells.Item(25,4).Value = 3
$ws.Cells.Item(25,5).Value = -33.333333333333
$ws.Cells.Item(25,6).Value = 24
$ws.Cells.Item(25,7).Value = 21
$ws.Cells.Item(25,8).Value = 14.285714285714
$ws.Cells.Item(25,9).Value = 255
$ws.Cells.Item(25,10).Value = 259
$ws.Cells.Item(25,11).Value = -1.544401544401
$ws.Cells.Item(25,12).Value = -18.789808917197

$ws.Cells.Item(26,3).Value = 22
$ws.Cells.Item(26,4).Value = 14
$ws.Cells.Item(26,5).Value = 57.142857142857
$ws.Cells.Item(26,6).Value = 90
$ws.Cells.Item(26,7).Value = 69
$ws.Cells.Item(26,8).Value = 30.434782608695
$ws.Cells.Item(26,9).Value = 789
$ws.Cells.Item(26,10).Value = 659
$ws.Cells.Item(26,11).Value = 19.726858877086
$ws.Cells.Item(26,12).Value = 25.837320574162
$ws.Cells.Item(26,13).Value = 17.062314540059

$ws.Cells.Item(27,6).Value = 8
$ws.Cells.Item(27,7).Value = 3
$ws.Cells.Item(27,8).Value = 166.666666666667
$ws.Cells.Item(27,9).Value = 45
$ws.Cells.Item(27,11).Value = 4.651162790697
$ws.Cells.Item(27,12).Value = -23.728813559322

$ws.Cells.Item(28,3).Value = 3
$ws.Cells.Item(28,4).Value = 1
$ws.Cells.Item(28,5).Value = 200
$ws.Cells.Item(28,6).Value = 4
$ws.Cells.Item(28,8).Value = -42.857142857142
$ws.Cells.Item(28,9).Value = 61
$ws.Cells.Item(28,10).Value = 54
$ws.Cells.Item(28,11).Value = 12.962962962963
$ws.Cells.Item(28,12).Value = 22

$ws.Cells.Item(29,13).Value = -53.061224489795
$ws.Cells.Item(29,14).Value = -77

$ws.Cells.Item(30,13).Value = -62.5
$ws.Cells.Item(30,14).Value = -84.210526315789
